$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "27.583.13"
$ws.Cells.Item(2, 5).Value = "  -1.09%  "

$ws.Cells.Item(3, 4).Value = "1.596.28"
$ws.Cells.Item(3, 5).Value = "  -2.08%  "

$ws.Cells.Item(4, 5).Value = "  +0.43%  "

$ws.Cells.Item(5, 5).Value = "  -1.51%  "

$ws.Cells.Item(6, 4).Value = "'0.502"
$ws.Cells.Item(6, 5).Value = "  -3.75%  "

$ws.Cells.Item(7, 5).Value = "  +0.47%  "

$ws.Cells.Item(8, 4).Value = "'22.36"
$ws.Cells.Item(8, 5).Value = "  -4.34%  "

$ws.Cells.Item(9, 5).Value = "  -1.98%  "

$ws.Cells.Item(10, 5).Value = "  -3.28%  "

$ws.Cells.Item(11, 5).Value = "  -1.75%  "

$ws.Cells.Item(12, 4).Value = "1.822.80"
$ws.Cells.Item(12, 5).Value = "  -2.07%  "

$ws.Cells.Item(13, 4).Value = "1.604.65"
$ws.Cells.Item(13, 5).Value = "  -1.56%  "

$ws.Cells.Item(14, 5).Value = "  -4.11%  "

$ws.Cells.Item(15, 4).Value = "'0.536"
$ws.Cells.Item(15, 5).Value = "  -4.52%  "

$ws.Cells.Item(16, 4).Value = "'63.39"
$ws.Cells.Item(16, 5).Value = "  -2.94%  "

$ws.Cells.Item(17, 4).Value = "27.588.13"
$ws.Cells.Item(17, 5).Value = "  -1.07%  "

$ws.Cells.Item(18, 4).Value = "'217.46"
$ws.Cells.Item(18, 5).Value = "  -5.17%  "

$ws.Cells.Item(19, 5).Value = "  -4.19%  "

$ws.Cells.Item(20, 4).Value = "0.0₃0694"
$ws.Cells.Item(20, 5).Value = "  -3.63%  "

$ws.Cells.Item(21, 5).Value = "  +0.53%  "

$ws.Cells.Item(22, 5).Value = "  -3.83%  "

$ws.Cells.Item(23, 5).Value = "  -4.51%  "

$ws.Cells.Item(24, 5).Value = "  -3.19%  "

$ws.Cells.Item(25, 4).Value = "'152.99"
$ws.Cells.Item(25, 5).Value = "  -0.70%  "

$ws.Cells.Item(26, 2).Value = "Cosmos"
$ws.Cells.Item(26, 3).Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Cells.Item(26, 4).Value = "'6.76"
$ws.Cells.Item(26, 5).Value = "  -1.96%  "

$ws.Cells.Item(27, 2).Value = "BinanceUSD"
$ws.Cells.Item(27, 3).Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Cells.Item(27, 4).Value = "'1.00"
$ws.Cells.Item(27, 5).Value = "  +0.47%  "

$ws.Cells.Item(28, 4).Value = "'15.10"
$ws.Cells.Item(28, 5).Value = "  -2.62%  "

$ws.Cells.Item(29, 5).Value = "  -3.89%  "

$ws.Cells.Item(30, 5).Value = "  -1.40%  "

$ws.Cells.Item(31, 5).Value = "  -3.09%  "

$ws.Cells.Item(32, 4).Value = "'3.26"
$ws.Cells.Item(32, 5).Value = "  -4.84%  "

$ws.Cells.Item(33, 4).Value = "1.370.95"
$ws.Cells.Item(33, 5).Value = "  -1.60%  "

$ws.Cells.Item(34, 4).Value = "'2.96"
$ws.Cells.Item(34, 5).Value = "  -5.10%  "

$ws.Cells.Item(35, 5).Value = "  -3.88%  "

$ws.Cells.Item(36, 4).Value = "'0.968"
$ws.Cells.Item(36, 5).Value = "  -5.10%  "

$ws.Cells.Item(37, 5).Value = "  -1.14%  "

$ws.Cells.Item(38, 5).Value = "  -3.20%  "

$ws.Cells.Item(39, 4).Value = "'0.541"
$ws.Cells.Item(39, 5).Value = "  -2.96%  "

$ws.Cells.Item(40, 4).Value = "'0.813"
$ws.Cells.Item(40, 5).Value = "  -4.53%  "

$ws.Cells.Item(41, 5).Value = "  +0.46%  "

$ws.Cells.Item(42, 4).Value = "'0.975"
$ws.Cells.Item(42, 5).Value = "  -4.05%  "

$ws.Cells.Item(43, 5).Value = "  -2.56%  "

$ws.Cells.Item(44, 5).Value = "  -1.22%  "

$ws.Cells.Item(45, 4).Value = "'64.06"
$ws.Cells.Item(45, 5).Value = "  -2.43%  "

$ws.Cells.Item(46, 4).Value = "'2.17"
$ws.Cells.Item(46, 5).Value = "  +0.94%  "

$ws.Cells.Item(47, 4).Value = "1.733.62"
$ws.Cells.Item(47, 5).Value = "  -2.13%  "

$ws.Cells.Item(48, 4).Value = "'87.30"
$ws.Cells.Item(48, 5).Value = "  -1.10%  "

$ws.Cells.Item(49, 5).Value = "  -3.37%  "

$ws.Cells.Item(50, 4).Value = "'0.0971"
$ws.Cells.Item(50, 5).Value = "  -4.27%  "

$ws.Cells.Item(51, 5).Value = "  -1.31%  "
